$wb = $excel.ActiveWorkbook

# Rename the 'Codelists' sheet to 'Cells'
$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# Make the 'Cells' sheet the active sheet/tab
$ws.Activate()

# Set the selection on the 'Cells' sheet to G13
$ws.Range("G13").Select()
